$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (names, URLs) - safe to assign directly
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('B20').Value = 'MCDex'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('B23').Value = 'UpBots'
$ws.Range('C23').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'

# Numeric/percent-looking text columns - must preserve as literal text
# Use a temporary formula producing the exact string, then convert formula to static value
# via Copy + PasteSpecial(xlPasteValues) so Excel does not auto-convert the text to a number
# and does not alter the cell style/number format.
$ws.Range('D2').Formula = "=""309.55"""
$ws.Range('D2').Copy()
$ws.Range('D2').PasteSpecial(-4163)
$ws.Range('E2').Formula = "=""0.37%"""
$ws.Range('E2').Copy()
$ws.Range('E2').PasteSpecial(-4163)
$ws.Range('D3').Formula = "=""37.12"""
$ws.Range('D3').Copy()
$ws.Range('D3').PasteSpecial(-4163)
$ws.Range('E3').Formula = "=""-1.83%"""
$ws.Range('E3').Copy()
$ws.Range('E3').PasteSpecial(-4163)
$ws.Range('D4').Formula = "=""5.126"""
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Formula = "=""1.70%"""
$ws.Range('E4').Copy()
$ws.Range('E4').PasteSpecial(-4163)
$ws.Range('D5').Formula = "=""0.07768"""
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Formula = "=""-1.60%"""
$ws.Range('E5').Copy()
$ws.Range('E5').PasteSpecial(-4163)
$ws.Range('D6').Formula = "=""4.393"""
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Formula = "=""0.04%"""
$ws.Range('E6').Copy()
$ws.Range('E6').PasteSpecial(-4163)
$ws.Range('D7').Formula = "=""8.205"""
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Formula = "=""-0.32%"""
$ws.Range('E7').Copy()
$ws.Range('E7').PasteSpecial(-4163)
$ws.Range('D8').Formula = "=""1.877"""
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Formula = "=""-8.54%"""
$ws.Range('E8').Copy()
$ws.Range('E8').PasteSpecial(-4163)
$ws.Range('D9').Formula = "=""0.9188"""
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Formula = "=""-0.90%"""
$ws.Range('E9').Copy()
$ws.Range('E9').PasteSpecial(-4163)
$ws.Range('D10').Formula = "=""0.1190"""
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Formula = "=""-7.30%"""
$ws.Range('E10').Copy()
$ws.Range('E10').PasteSpecial(-4163)
$ws.Range('D11').Formula = "=""0.1889"""
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Formula = "=""-1.02%"""
$ws.Range('E11').Copy()
$ws.Range('E11').PasteSpecial(-4163)
$ws.Range('D12').Formula = "=""0.09171"""
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Formula = "=""5.21%"""
$ws.Range('E12').Copy()
$ws.Range('E12').PasteSpecial(-4163)
$ws.Range('D13').Formula = "=""0.03415"""
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Formula = "=""-1.16%"""
$ws.Range('E13').Copy()
$ws.Range('E13').PasteSpecial(-4163)
$ws.Range('D14').Formula = "=""0.09683"""
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Formula = "=""-0.67%"""
$ws.Range('E14').Copy()
$ws.Range('E14').PasteSpecial(-4163)
$ws.Range('D15').Formula = "=""0.001375"""
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Formula = "=""-1.75%"""
$ws.Range('E15').Copy()
$ws.Range('E15').PasteSpecial(-4163)
$ws.Range('D16').Formula = "=""0.005801"""
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Formula = "=""-5.39%"""
$ws.Range('E16').Copy()
$ws.Range('E16').PasteSpecial(-4163)
$ws.Range('D17').Formula = "=""3.551"""
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Formula = "=""0.02%"""
$ws.Range('E17').Copy()
$ws.Range('E17').PasteSpecial(-4163)
$ws.Range('D18').Formula = "=""3.054"""
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Formula = "=""-1.02%"""
$ws.Range('E18').Copy()
$ws.Range('E18').PasteSpecial(-4163)
$ws.Range('D19').Formula = "=""0.3398"""
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Formula = "=""-1.27%"""
$ws.Range('E19').Copy()
$ws.Range('E19').PasteSpecial(-4163)
$ws.Range('D20').Formula = "=""5.271"""
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Formula = "=""5.11%"""
$ws.Range('E20').Copy()
$ws.Range('E20').PasteSpecial(-4163)
$ws.Range('D21').Formula = "=""0.1275"""
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Formula = "=""-2.05%"""
$ws.Range('E21').Copy()
$ws.Range('E21').PasteSpecial(-4163)
$ws.Range('D22').Formula = "=""0.2594"""
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Formula = "=""2.92%"""
$ws.Range('E22').Copy()
$ws.Range('E22').PasteSpecial(-4163)
$ws.Range('D23').Formula = "=""0.02107"""
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Formula = "=""5,595.56%"""
$ws.Range('E23').Copy()
$ws.Range('E23').PasteSpecial(-4163)
$ws.Range('D24').Formula = "=""0.04337"""
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Formula = "=""0.17%"""
$ws.Range('E24').Copy()
$ws.Range('E24').PasteSpecial(-4163)
$ws.Range('D25').Formula = "=""0.001200"""
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Formula = "=""-1.66%"""
$ws.Range('E25').Copy()
$ws.Range('E25').PasteSpecial(-4163)
$ws.Range('D26').Formula = "=""0.004246"""
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Formula = "=""-7.67%"""
$ws.Range('E26').Copy()
$ws.Range('E26').PasteSpecial(-4163)
$ws.Range('D27').Formula = "=""0.0001312"""
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Formula = "=""-63.49%"""
$ws.Range('E27').Copy()
$ws.Range('E27').PasteSpecial(-4163)
$ws.Range('D39').Formula = "=""0.02065"""
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Formula = "=""-8.59%"""
$ws.Range('E39').Copy()
$ws.Range('E39').PasteSpecial(-4163)
$ws.Range('D40').Formula = "=""0.05026"""
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Formula = "=""-0.64%"""
$ws.Range('E40').Copy()
$ws.Range('E40').PasteSpecial(-4163)
$ws.Range('D41').Formula = "=""0.007667"""
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Formula = "=""2.23%"""
$ws.Range('E41').Copy()
$ws.Range('E41').PasteSpecial(-4163)
$ws.Range('D42').Formula = "=""0.009806"""
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Formula = "=""-1.31%"""
$ws.Range('E42').Copy()
$ws.Range('E42').PasteSpecial(-4163)
$ws.Range('D43').Formula = "=""0.1344"""
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Formula = "=""-1.20%"""
$ws.Range('E43').Copy()
$ws.Range('E43').PasteSpecial(-4163)
$ws.Range('D44').Formula = "=""0.002173"""
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Formula = "=""3.70%"""
$ws.Range('E44').Copy()
$ws.Range('E44').PasteSpecial(-4163)
$ws.Range('D45').Formula = "=""0.008764"""
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Formula = "=""-0.74%"""
$ws.Range('E45').Copy()
$ws.Range('E45').PasteSpecial(-4163)
$ws.Range('D46').Formula = "=""0.00006711"""
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Formula = "=""3.18%"""
$ws.Range('E46').Copy()
$ws.Range('E46').PasteSpecial(-4163)
$ws.Range('E47').Formula = "=""-0.15%"""
$ws.Range('E47').Copy()
$ws.Range('E47').PasteSpecial(-4163)
$ws.Range('D48').Formula = "=""0.001202"""
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Formula = "=""-0.14%"""
$ws.Range('E48').Copy()
$ws.Range('E48').PasteSpecial(-4163)
$ws.Range('D49').Formula = "=""0.002938"""
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Formula = "=""-2.24%"""
$ws.Range('E49').Copy()
$ws.Range('E49').PasteSpecial(-4163)
$ws.Range('E50').Formula = "=""-0.15%"""
$ws.Range('E50').Copy()
$ws.Range('E50').PasteSpecial(-4163)
$ws.Range('E51').Formula = "=""-0.15%"""
$ws.Range('E51').Copy()
$ws.Range('E51').PasteSpecial(-4163)
$excel.CutCopyMode = 0
